{"js": "// Merge the multi-run \"word-by-word\" paragraphs (Title, Author, Abstract)\n// into single runs containing the full sentence text, leaving all other\n// content (formatting, styles, other paragraphs) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,style,text\");\nawait context.sync();\n\nconst replacements = {\n  \"Title\": \"Questions: Laws of indices\",\n  \"Author\": \"Isabella Lewis, Akshat Srivastava\",\n  \"Abstract\": \"A selection of questions for the study guide on laws of indices.\"\n};\n\nfor (const para of paragraphs.items) {\n  const style = para.style;\n  if (Object.prototype.hasOwnProperty.call(replacements, style)) {\n    const target = replacements[style];\n    // Always rewrite: the paragraph's aggregate text may already equal the\n    // target even though it is still split across many single-word runs in\n    // the underlying OOXML. Replacing collapses it into one run.\n    para.getRange().insertText(target, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Merge the multi-run \"word-by-word\" paragraphs (Title, Author, Abstract)\n# into single runs containing the full sentence text, leaving all other\n# content (formatting, styles, other paragraphs) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"Title\"    = \"Questions: Laws of indices\"\n    \"Author\"   = \"Isabella Lewis, Akshat Srivastava\"\n    \"Abstract\" = \"A selection of questions for the study guide on laws of indices.\"\n}\n\nforeach ($para in $d.Paragraphs) {\n    $styleName = $para.Style.NameLocal\n    if ($replacements.ContainsKey($styleName)) {\n        # Re-derive a fresh document-level range (rather than reusing the\n        # paragraph's own .Range object) covering the paragraph's full text\n        # span, so the assignment replaces ALL runs in the paragraph with a\n        # single new run, instead of only overwriting the first run's text\n        # in place (which is what happens if `$para.Range.Text = ...` is\n        # used directly on the paragraph's own cached Range object).\n        $start = $para.Range.Start\n        $end = $para.Range.End\n        $r = $d.Range($start, $end)\n        $r.Text = $replacements[$styleName]\n    }\n}\n"}
